# Refresh the customers test data: rename the two sample customers and
# change their pin/email values so the Guru99 "add existing customer"
# test doesn't collide with previously-created accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pin column first
$ws.Range("G2").Value = '"222222"'
$ws.Range("G3").Value = '"333333"'

# Then the customer names
$ws.Range("A2").Value = "alexiss"
$ws.Range("A3").Value = "Dinas"

# Then the email addresses
$ws.Range("I2").Value = "artahAlsd+1@gmail.com"
$ws.Range("I3").Value = "stsADDdDN+Din@gmail.com"

# Column width tweaks (mobile/email columns)
$ws.Range("H1").EntireColumn.ColumnWidth = 15.6
$ws.Range("I1").EntireColumn.ColumnWidth = 27.8

# Update the active selection
$ws.Range("D7").Select()
